$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.295069333333332
$ws.Range("H2").Value = 27.885208
$ws.Range("I2").Value = 0.2851098797714356
$ws.Range("J2").Value = 0.2851098797714357
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 1087.937328174701
$ws.Range("R2").Value = 9791.435953572312
$ws.Range("S2").Value = 0.09252858089822719
$ws.Range("T2").Value = 0.0925285808982272

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.295069333333332
$ws.Range("H3").Value = 27.885208
$ws.Range("I3").Value = 0.2851098797714356
$ws.Range("J3").Value = 0.2851098797714357
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 944.1934898959217
$ws.Range("R3").Value = 8497.741409063296
$ws.Range("S3").Value = 0.08030323204369837
$ws.Range("T3").Value = 0.08030323204369837

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.295069333333332
$ws.Range("H4").Value = 27.885208
$ws.Range("I4").Value = 0.2851098797714356
$ws.Range("J4").Value = 0.2851098797714357
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 1320.14885404406
$ws.Range("R4").Value = 11881.33968639654
$ws.Range("S4").Value = 0.1122780668295101
$ws.Range("T4").Value = 0.1122780668295101

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.588082333333332
$ws.Range("H5").Value = 25.764247
$ws.Range("I5").Value = 0.2634242988100204
$ws.Range("J5").Value = 0.2634242988100204
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 1005.18834371302
$ws.Range("R5").Value = 9046.695093417182
$ws.Range("S5").Value = 0.08549081695289515
$ws.Range("T5").Value = 0.08549081695289515

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.588082333333332
$ws.Range("H6").Value = 25.764247
$ws.Range("I6").Value = 0.2634242988100204
$ws.Range("J6").Value = 0.2634242988100204
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 872.3777240417404
$ws.Range("R6").Value = 7851.399516375664
$ws.Range("S6").Value = 0.07419533342810854
$ws.Range("T6").Value = 0.07419533342810852

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.588082333333332
$ws.Range("H7").Value = 25.764247
$ws.Range("I7").Value = 0.2634242988100204
$ws.Range("J7").Value = 0.2634242988100204
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 1219.737760333655
$ws.Range("R7").Value = 10977.63984300289
$ws.Range("S7").Value = 0.1037381484290168
$ws.Range("T7").Value = 0.1037381484290168

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.71855733333333
$ws.Range("H8").Value = 44.155672
$ws.Range("I8").Value = 0.4514658214185439
$ws.Range("J8").Value = 0.4514658214185439
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 1722.727111070446
$ws.Range("R8").Value = 15504.54399963401
$ws.Range("S8").Value = 0.1465171666916591
$ws.Range("T8").Value = 0.1465171666916591

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.71855733333333
$ws.Range("H9").Value = 44.155672
$ws.Range("I9").Value = 0.4514658214185439
$ws.Range("J9").Value = 0.4514658214185439
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 1495.111603412807
$ws.Range("R9").Value = 13456.00443071527
$ws.Range("S9").Value = 0.1271585700440691
$ws.Range("T9").Value = 0.1271585700440691

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.71855733333333
$ws.Range("H10").Value = 44.155672
$ws.Range("I10").Value = 0.4514658214185439
$ws.Range("J10").Value = 0.4514658214185439
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 2090.429441671922
$ws.Range("R10").Value = 18813.8649750473
$ws.Range("S10").Value = 0.1777900846828157
$ws.Range("T10").Value = 0.1777900846828157
